$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The NATMI "Resolving-Mac" cluster was dropped as a possible *target* cluster
# for this ligand-receptor pair, so the four rows describing edges that ended
# in "Resolving-Mac" (previously rows 14-17, i.e. every remaining sending
# cluster paired with D = "Resolving-Mac") are removed entirely.
$ws.Range("A14:T17").EntireRow.Delete()

# The remaining 12 rows (every combination of the 4 sending clusters against
# the 3 surviving target clusters: ECs, FAPs, MuSCs) keep their row position,
# but the Sending cluster / Target cluster labels are rewritten and every
# NATMI-derived metric (columns G through T) is recomputed against the new
# TPM values.
$data = @(
    @("ECs",           "ECs",   4.232924,          12.698772,   0.05792409824508498, 0.05792409824508497, 2, 0.6666666666666666, 0.227228,           0.681684,  0.2376267857721762, 0.2376267857721762, 0.9618388546719999, 8.656549692047999, 0.0137643172847313,   0.01376431728473129),
    @("ECs",           "FAPs",  4.232924,          12.698772,   0.05792409824508498, 0.05792409824508497, 3, 1,                  0.6537306666666667, 1.961192,  0.6836477770376096, 0.6836477770376095, 2.767192228469333,  24.904730056224,   0.03959968100216045,  0.03959968100216044),
    @("ECs",           "MuSCs", 4.232924,          12.698772,   0.05792409824508498, 0.05792409824508497, 1, 0.3333333333333333, 0.07528033333333332,0.225841,  0.0787254371902143, 0.0787254371902143, 0.3186559296946666, 2.867903367252,    0.00456009995819324,  0.004560099958193239),
    @("FAPs",          "ECs",   23.77965533333333, 71.338966,   0.3254051080913003, 0.3254051080913002,  2, 0.6666666666666666, 0.227228,           0.681684,  0.2376267857721762, 0.2376267857721762, 5.403403522082667,  48.630631698744,   0.07732496990958325,  0.07732496990958324),
    @("FAPs",          "FAPs",  23.77965533333333, 71.338966,   0.3254051080913003, 0.3254051080913002,  3, 1,                  0.6537306666666667, 1.961192,  0.6836477770376096, 0.6836477770376095, 15.54548993416356,  139.909409407472,  0.2224624787833005,   0.2224624787833004),
    @("FAPs",          "MuSCs", 23.77965533333333, 71.338966,   0.3254051080913003, 0.3254051080913002,  1, 0.3333333333333333, 0.07528033333333332,0.225841,  0.0787254371902143, 0.0787254371902143, 1.790140380045111,  16.111263420406,   0.02561765939841656,  0.02561765939841655),
    @("MuSCs",         "ECs",   7.006365333333332, 21.019096,   0.09587637148905993, 0.09587637148905992,2, 0.6666666666666666, 0.227228,           0.681684,  0.2376267857721762, 0.2376267857721762, 1.592042381962666,  14.328381437664,   0.02278279398844443,  0.02278279398844442),
    @("MuSCs",         "FAPs",  7.006365333333332, 21.019096,   0.09587637148905993, 0.09587637148905992,3, 1,                  0.6537306666666667, 1.961192,  0.6836477770376096, 0.6836477770376095, 4.580275880270221,  41.22248292243199, 0.06554566823892788,  0.06554566823892785),
    @("MuSCs",         "MuSCs", 7.006365333333332, 21.019096,   0.09587637148905993, 0.09587637148905992,1, 0.3333333333333333, 0.07528033333333332,0.225841,  0.0787254371902143, 0.0787254371902143, 0.5274415177484443, 4.746973659735999, 0.007547909261687641, 0.00754790926168764),
    @("Resolving-Mac", "ECs",   38.058136,         114.174408,  0.5207944221745548, 0.5207944221745548,  2, 0.6666666666666666, 0.227228,           0.681684,  0.2376267857721762, 0.2376267857721762, 8.647874127007999,  77.830867143072,   0.1237547045894172,   0.1237547045894172),
    @("Resolving-Mac", "FAPs",  38.058136,         114.174408,  0.5207944221745548, 0.5207944221745548,  3, 1,                  0.6537306666666667, 1.961192,  0.6836477770376096, 0.6836477770376095, 24.87977061937067,  223.917935574336,  0.3560399490132208,   0.3560399490132207),
    @("Resolving-Mac", "MuSCs", 38.058136,         114.174408,  0.5207944221745548, 0.5207944221745548,  1, 0.3333333333333333, 0.07528033333333332,0.225841,  0.0787254371902143, 0.0787254371902143, 2.865029164125333,  25.785262477128,   0.04099976857191687,  0.04099976857191687)
)

$columns = @("A", "D", "G", "H", "I", "J", "K", "L", "M", "N", "O", "P", "Q", "R", "S", "T")

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $columns.Length; $j++) {
        $ws.Range("$($columns[$j])$rowNum").Value = $values[$j]
    }
}
